$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2016806722689076
$ws.Range("C2").Value = 0.5294117647058824
$ws.Range("J2").Value = 0.02941176470588235
$ws.Range("P2").Value = 0.1512605042016807
$ws.Range("S2").Value = 0.08823529411764706
$ws.Range("B3").Value = 0.01526717557251908
$ws.Range("C3").Value = 0.03816793893129771
$ws.Range("J3").Value = 0.05343511450381679
$ws.Range("P3").Value = 0.7251908396946565
$ws.Range("S3").Value = 0.1679389312977099
$ws.Range("J4").Value = 0.05405405405405406
$ws.Range("O4").Value = 0.02702702702702703
$ws.Range("P4").Value = 0.7297297297297297
$ws.Range("S4").Value = 0.1891891891891892
$ws.Range("B6").Value = 0.1014492753623188
$ws.Range("D6").Value = 0.01449275362318841
$ws.Range("F6").Value = 0.05797101449275362
$ws.Range("J6").Value = 0.2173913043478261
$ws.Range("O6").Value = 0.03623188405797102
$ws.Range("Q6").Value = 0.1521739130434783
$ws.Range("R6").Value = 0.06521739130434782
$ws.Range("S6").Value = 0.355072463768116
$ws.Range("B7").Value = 0.1363636363636364
$ws.Range("D7").Value = 0.03636363636363636
$ws.Range("F7").Value = 0.07272727272727272
$ws.Range("J7").Value = 0.1272727272727273
$ws.Range("Q7").Value = 0.1727272727272727
$ws.Range("R7").Value = 0.07272727272727272
$ws.Range("S7").Value = 0.3818181818181818
$ws.Range("B8").Value = 0.1111111111111111
$ws.Range("D8").Value = 0.01388888888888889
$ws.Range("F8").Value = 0.05208333333333334
$ws.Range("J8").Value = 0.1354166666666667
$ws.Range("O8").Value = 0.003472222222222222
$ws.Range("Q8").Value = 0.1319444444444444
$ws.Range("R8").Value = 0.09722222222222222
$ws.Range("S8").Value = 0.4548611111111111
$ws.Range("B9").Value = 0.08552631578947369
$ws.Range("D9").Value = 0.04605263157894737
$ws.Range("E9").Value = 0.006578947368421052
$ws.Range("F9").Value = 0.06578947368421052
$ws.Range("J9").Value = 0.09210526315789473
$ws.Range("O9").Value = 0.0131578947368421
$ws.Range("Q9").Value = 0.1842105263157895
$ws.Range("R9").Value = 0.03947368421052631
$ws.Range("S9").Value = 0.4671052631578947
$ws.Range("B10").Value = 0.1267281105990783
$ws.Range("D10").Value = 0.02304147465437788
$ws.Range("F10").Value = 0.06682027649769585
$ws.Range("J10").Value = 0.1589861751152074
$ws.Range("O10").Value = 0.0195852534562212
$ws.Range("Q10").Value = 0.1658986175115207
$ws.Range("R10").Value = 0.07373271889400922
$ws.Range("S10").Value = 0.3652073732718894
$ws.Range("F11").Value = 0.005434782608695652
$ws.Range("G11").Value = 0.1630434782608696
$ws.Range("J11").Value = 0.07608695652173914
$ws.Range("K11").Value = 0.2119565217391304
$ws.Range("L11").Value = 0.5326086956521739
$ws.Range("S11").Value = 0.0108695652173913
$ws.Range("G12").Value = 0.7128712871287128
$ws.Range("J12").Value = 0.2079207920792079
$ws.Range("K12").Value = 0.009900990099009901
$ws.Range("L12").Value = 0.0396039603960396
$ws.Range("S12").Value = 0.0297029702970297
$ws.Range("G13").Value = 0.65
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.1
$ws.Range("F15").Value = 0.01273885350318471
$ws.Range("H15").Value = 0.1592356687898089
$ws.Range("I15").Value = 0.1082802547770701
$ws.Range("J15").Value = 0.3375796178343949
$ws.Range("K15").Value = 0.07006369426751592
$ws.Range("M15").Value = 0.006369426751592357
$ws.Range("O15").Value = 0.06369426751592357
$ws.Range("S15").Value = 0.2420382165605096
$ws.Range("F16").Value = 0.006451612903225806
$ws.Range("H16").Value = 0.1935483870967742
$ws.Range("I16").Value = 0.07741935483870968
$ws.Range("J16").Value = 0.4387096774193548
$ws.Range("K16").Value = 0.09677419354838709
$ws.Range("M16").Value = 0.03225806451612903
$ws.Range("O16").Value = 0.03225806451612903
$ws.Range("S16").Value = 0.1225806451612903
$ws.Range("F17").Value = 0.01606425702811245
$ws.Range("H17").Value = 0.2208835341365462
$ws.Range("I17").Value = 0.08835341365461848
$ws.Range("J17").Value = 0.3253012048192771
$ws.Range("K17").Value = 0.08433734939759036
$ws.Range("M17").Value = 0.01606425702811245
$ws.Range("O17").Value = 0.04819277108433735
$ws.Range("S17").Value = 0.2008032128514056
$ws.Range("F18").Value = 0.04504504504504504
$ws.Range("H18").Value = 0.1621621621621622
$ws.Range("I18").Value = 0.09009009009009009
$ws.Range("J18").Value = 0.4324324324324325
$ws.Range("K18").Value = 0.08108108108108109
$ws.Range("O18").Value = 0.04504504504504504
$ws.Range("S18").Value = 0.1441441441441441
$ws.Range("F19").Value = 0.0121012101210121
$ws.Range("H19").Value = 0.1804180418041804
$ws.Range("I19").Value = 0.1012101210121012
$ws.Range("J19").Value = 0.3718371837183718
$ws.Range("K19").Value = 0.09680968096809681
$ws.Range("M19").Value = 0.0143014301430143
$ws.Range("O19").Value = 0.08030803080308031
$ws.Range("S19").Value = 0.1507731958762887
